$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new text looks like a plain number; force Text format
# so Excel keeps them as literal strings (matching the source "inlineStr" cells)
# instead of silently converting to a numeric value (which would also silently
# drop meaningful trailing zeros, e.g. "146.90" -> 146.9).
$textCells = @("D4", "D5", "D6", "D7", "D13", "D14", "D19", "D20", "D21", "D22", "D24", "D25", "D27", "D28", "D30", "D35", "D36", "D37", "D38", "D41", "D42", "D44", "D45", "D46", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "61.807.89"
$ws.Range("E2").Value = "  +2.50%  "
$ws.Range("D3").Value = "2.385.14"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "553.59"
$ws.Range("E5").Value = "  +2.44%  "
$ws.Range("D6").Value = "141.51"
$ws.Range("E6").Value = "  +3.82%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("D9").Value = "2.383.78"
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("E10").Value = "  +5.09%  "
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").Value = "0.354"
$ws.Range("E13").Value = "  +4.53%  "
$ws.Range("D14").Value = "25.78"
$ws.Range("E14").Value = "  +6.09%  "
$ws.Range("E15").Value = "  +10.58%  "
$ws.Range("D16").Value = "2.813.55"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").Value = "61.744.00"
$ws.Range("E17").Value = "  +2.29%  "
$ws.Range("D18").Value = "2.380.41"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "11.04"
$ws.Range("E19").Value = "  +5.37%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "4.19"
$ws.Range("E20").Value = "  +3.63%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "324.16"
$ws.Range("E21").Value = "  +4.57%  "
$ws.Range("D22").Value = "6.72"
$ws.Range("E22").Value = "  +2.77%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "1.79"
$ws.Range("E24").Value = "  -3.84%  "
$ws.Range("D25").Value = "64.47"
$ws.Range("E25").Value = "  +2.61%  "
$ws.Range("E26").Value = "  +5.39%  "
$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").Value = "539.10"
$ws.Range("E27").Value = "  +9.37%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").Value = "2.501.70"
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("D30").Value = "8.31"
$ws.Range("E30").Value = "  +5.22%  "
$ws.Range("D31").Value = "0.0₃0924"
$ws.Range("E31").Value = "  +5.61%  "
$ws.Range("E32").Value = "  +3.84%  "
$ws.Range("E33").Value = "  +4.12%  "
$ws.Range("E34").Value = "  +4.55%  "
$ws.Range("D35").Value = "1.54"
$ws.Range("E35").Value = "  +2.14%  "
$ws.Range("D36").Value = "5.76"
$ws.Range("E36").Value = "  +11.04%  "
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "4.78"
$ws.Range("E38").Value = "  +5.35%  "
$ws.Range("E39").Value = "  +9.07%  "
$ws.Range("E40").Value = "  +3.57%  "
$ws.Range("D41").Value = "18.61"
$ws.Range("E41").Value = "  +2.04%  "
$ws.Range("D42").Value = "146.90"
$ws.Range("E42").Value = "  +7.29%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "41.54"
$ws.Range("E44").Value = "  +4.09%  "
$ws.Range("D45").Value = "149.07"
$ws.Range("E45").Value = "  +5.73%  "
$ws.Range("D46").Value = "2.22"
$ws.Range("E46").Value = "  +6.37%  "
$ws.Range("E47").Value = "  +3.19%  "
$ws.Range("D48").Value = "0.0530"
$ws.Range("E48").Value = "  +4.77%  "
$ws.Range("D49").Value = "20.24"
$ws.Range("E49").Value = "  +4.99%  "
$ws.Range("E50").Value = "  +3.34%  "
$ws.Range("D51").Value = "0.0908"
$ws.Range("E51").Value = "  +1.76%  "
